# Update vendas_atipicas_atual.xlsx sheet with new "vendas atipicas" data
# Rows 2-23 are replaced with refreshed values and a new row 24 is appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array: Dia(A), quantidade_atipica(B - text date), C, cliente(D), id_produto(E - text),
# produto(F), estoque_atualizado(G), critico(H - boolean)
$data = @(
    @(22, "2025-03-13", 360, "AMMAC INDUSTRIA E COMERCIO DE ALIMENTOS LTDA", "000029", "ESPONJA MULTIUSO JEITOSA", 2416, 0),
    @(19, "2025-03-13", 50, "E A H EMPRESA AMAZONENSE DE HOTELARIA LTDA", "000440", "LUVA DE LIMPEZA LATEX TAM M AMARELA TOP REFOR NOBRE CA 47259", 292, 1),
    @(16, "2025-03-14", 27, "AMAZONIA REFEICOES E SERVICOS LTDA", "000152", "COPO DESCARTAVEL BRANCO CRISTALCOPO 180ML CX C\25", 56, 0),
    @(9, "2025-03-14", 11, "AMAZONIA REFEICOES E SERVICOS LTDA", "000924", "COPO POTE DESCARTAVEL TRANSP 100ML CX/20", 1, 0),
    @(13, "2025-03-17", 180, "MUSASHI DA AMAZONIA LTDA", "000415", "DETERGENTE DESENGRAX MAX PINE AUDAX 5L", 3, 1),
    @(21, "2025-03-17", 432, "DSB COMERCIO DE ALIMENTOS EIRELI", "000717", "DETERGENTE NEUTRO BRINORT 500ML", 351, 0),
    @(0, "2025-03-17", 40, "POTENCIAL HUMANO RECRUTAMENTO E SELECAO LTDA", "000078", "SACO DE LIXO 100L PRETO COMUM - PCT C/100 UND", 72, 0),
    @(3, "2025-03-17", 200, "MUSASHI DA AMAZONIA LTDA", "000243", "FIBRA DE LIMPEZA PESADA 102X260mm", 661, 0),
    @(6, "2025-03-17", 36, "POTENCIAL HUMANO RECRUTAMENTO E SELECAO LTDA", "000071", "SACO DE LIXO 200L PRETO COMUM - PCT C/100 UND", 7, 0),
    @(11, "2025-03-17", 400, "MUSASHI DA AMAZONIA LTDA", "000842", "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND", 426, 0),
    @(15, "2025-03-19", 25, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000276", "INSETICIDA SBP AEROSSOL 380ML", 50, 0),
    @(17, "2025-03-19", 30, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000146", "DESINFETANTE BRINORT 2L LAVANDA", 42, 0),
    @(1, "2025-03-19", 200, "ZARAPLAST DA AMAZONIA LTDA", "000098", "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM", 311, 0),
    @(5, "2025-03-24", 150, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000088", "VASSOURA PIACAVA 20 FUROS", 2, 0),
    @(7, "2025-03-25", 40, "MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.", "000288", "TOUCA DESCARTAVEL TNT TALGE PCT C/ 100 UND", 400, 0),
    @(10, "2025-03-25", 15, "AMAZONIA REFEICOES E SERVICOS LTDA", "000924", "COPO POTE DESCARTAVEL TRANSP 100ML CX/20", 1, 0),
    @(4, "2025-03-25", 120, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000350", "DESODORISADOR LADY AEROSSOL 360 ML LAVANDA", 1017, 1),
    @(12, "2025-03-25", 96, "MM DA AMAZONIA INDUSTRIA E COMERCIO DE PLASTICOS LTDA.", "000032", "LIMPADOR VEJA MULTIUSO GOLD 500ML", 1382, 0),
    @(2, "2025-03-25", 300, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000098", "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM", 311, 0),
    @(18, "2025-03-25", 200, "MAP SERVICOS DE CONSERVACAO - EIRELI", "000486", "NAFTALINA 40G PCT RUBI", 310, 0),
    @(20, "2025-03-25", 96, "MANJAR SERVICOS GERAIS SA", "000583", "LIMPA ALUMINIO BRINORT 500ML", 221, 0),
    @(8, "2025-03-26", 10, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000425", "COADOR DE CAFE EG (EXTRA GRANDE)", 5, 0),
    @(14, "2025-03-26", 200, "JURUA ESTALEIROS E NAVEGACAO LTDA", "000122", "SABAO EM PO ALA LAVANDA ROUPAS 400G", 58, 0)
)

# Columns that must stay TEXT (avoid Excel's automatic date / number coercion).
# B holds dates written as plain text ("2025-03-13"), E holds zero-padded codes
# ("000029") -- both would otherwise be auto-converted by Excel's smart entry.
$ws.Range("B2:B24").NumberFormat = "@"
$ws.Range("E2:E24").NumberFormat = "@"

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]       # A - Dia
    $ws.Cells.Item($row, 2).Value = [string]$rec[1]  # B - quantidade_atipica (date text)
    $ws.Cells.Item($row, 3).Value = $rec[2]       # C
    $ws.Cells.Item($row, 4).Value = [string]$rec[3]  # D - cliente
    $ws.Cells.Item($row, 5).Value = [string]$rec[4]  # E - id_produto
    $ws.Cells.Item($row, 6).Value = [string]$rec[5]  # F - produto
    $ws.Cells.Item($row, 7).Value = $rec[6]       # G - estoque_atualizado
    $ws.Cells.Item($row, 8).Value = [bool]($rec[7] -ne 0)  # H - critico
}

# Row 24 is new: give column A the same style (border/center/top alignment) used by
# the rest of column A (copied from A2's format).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
